$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values after repulling / recalculating data
$updates = @{
    2  = -5
    6  = -6
    7  = -4
    9  = -7
    10 = -10
    12 = 6
    14 = -5
    16 = 1
    17 = -4
    19 = 0
    23 = -9
    24 = -3
    25 = -7
    26 = -6
    27 = 5
    30 = -7
    35 = 3
    37 = -1
    41 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
